$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 119765.1405966122
$ws.Range("E6").Value = 0.005109743552133421
$ws.Range("F6").Value = 0.2210963616989756
$ws.Range("G6").Value = -1.383155733590929
$ws.Range("H6").Value = 11.97557868494806
$ws.Range("D8").Value = 121559.4643682876
$ws.Range("E8").Value = -0.02206597877392203
$ws.Range("F8").Value = 0.2060234106399831
$ws.Range("G8").Value = -1.16871571259097
$ws.Range("H8").Value = 8.75347982884247
$ws.Range("D9").Value = 123286.2422882786
$ws.Range("E9").Value = -0.05094705112303217
$ws.Range("F9").Value = 0.318136351082018
$ws.Range("G9").Value = -1.603784022362666
$ws.Range("H9").Value = 10.6654282228842
$ws.Range("D10").Value = 124539.5362482381
$ws.Range("E10").Value = -0.0958308110981263
$ws.Range("F10").Value = 0.4280735799053507
$ws.Range("G10").Value = -1.884474311310691
$ws.Range("H10").Value = 9.768235697739948
$ws.Range("D11").Value = 126539.1800181963
$ws.Range("E11").Value = -0.169981061677895
$ws.Range("F11").Value = 0.7361749700893891
$ws.Range("G11").Value = -2.564224769720463
$ws.Range("H11").Value = 12.42327081754483
$ws.Range("D14").Value = 116792.8286770761
$ws.Range("E14").Value = 0.1167721660652652
$ws.Range("F14").Value = 0.1133600116321851
$ws.Range("G14").Value = -0.3128662589603944
$ws.Range("H14").Value = 5.654085880062441
$ws.Range("D15").Value = 116839.2840796073
$ws.Range("E15").Value = 0.1059050861407755
$ws.Range("F15").Value = 0.1235509118824574
$ws.Range("G15").Value = -0.07201929273239523
$ws.Range("H15").Value = 7.611158043261241
$ws.Range("D17").Value = 117628.1019040164
$ws.Range("E17").Value = 0.09150261896082844
$ws.Range("F17").Value = 0.1068964240361746
$ws.Range("G17").Value = -0.2448844007750582
$ws.Range("H17").Value = 5.095459527522543
$ws.Range("D18").Value = 117664.1911717433
$ws.Range("E18").Value = 0.06778743546482378
$ws.Range("F18").Value = 0.1338276161296763
$ws.Range("G18").Value = -0.05651437813418732
$ws.Range("H18").Value = 5.723602935046212
$ws.Range("D19").Value = 117677.0956656522
$ws.Range("E19").Value = 0.05365971831452972
$ws.Range("F19").Value = 0.1384872011011651
$ws.Range("G19").Value = -0.368392055979987
$ws.Range("H19").Value = 6.909798936527114
$ws.Range("D20").Value = 118695.9430765598
$ws.Range("E20").Value = 0.06142500222440995
$ws.Range("F20").Value = 0.137386212883784
$ws.Range("G20").Value = -0.2957987716406761
$ws.Range("H20").Value = 5.8085666925921

Write-Host "Updated calibration values after removing sub-USD5 noise point from extrapolation."
